$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ponds")

# Fix typo: "Amoun (Gal)" -> "Amount (Gal)" for the header in G3
$ws.Range("G3").Value = "Amount (Gal)"

# Update the active selection as recorded in the saved view state
$ws.Range("G23").Select()
